$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$boldText = "Meta description"
$restText = ": Cheeky Fruits 6 Deluxe is a fruit-themed slot game with high-quality graphics, exciting sound effects, and a bonus symbol. Play for free and activate bonus modes."

$metaPara.Range.Text = $boldText + $restText

$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)
$boldRange.Bold = 1

# Leave a (formatting-less) leading empty run in the paragraph, matching the
# pattern used throughout the rest of the document's body paragraphs.
$leadIn = $d.Range($metaStart, $metaStart)
$leadIn.InsertBefore("")

# ---------------------------------------------------------------------------
# 2) Remove the paragraph that duplicated the page title in bold
#    ("Play Cheeky Fruits 6 Deluxe Free - Exciting Fruit-themed Slot Game")
#    near the end of the document, and 3) replace the text of the final
#    (italic) paragraph with the new image-prompt text, while keeping its
#    run formatting intact.
# ---------------------------------------------------------------------------
$oldImageText = "Cheeky Fruits 6 Deluxe is a fruit-themed slot game with high-quality graphics, exciting sound effects, and a bonus symbol. Play for free and activate bonus modes."
$newImageText = 'Create a feature image that fits the game "Cheeky Fruits 6 Deluxe". The image should be in cartoon style and should feature a happy Maya warrior with glasses.'
$oldBoldTitle = "Play Cheeky Fruits 6 Deluxe Free - Exciting Fruit-themed Slot Game"

$count = $d.Paragraphs.Count
$secondLastPara = $d.Paragraphs($count - 1)

if ($secondLastPara.Range.Text.TrimEnd() -eq $oldBoldTitle) {
    $secondLastPara.Range.Delete()
}

# Re-fetch the last paragraph (now at $count - 1) and replace its text in
# place, preserving the leading empty run and the italic run formatting.
# Use Find purely to locate the exact sub-range (no in-place Find replace,
# since that would mangle straight quotes into curly ones); then assign
# .Text directly on that located sub-range so the existing run (and its
# <w:i/> formatting) is reused instead of being rebuilt.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$textRange = $lastPara.Range.Duplicate()
$found = $textRange.Find.Execute($oldImageText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $textRange.Text = $newImageText
}
